# Update the cryptos list with freshly scraped prices / 1h volume deltas.
#
# Price-column (D) values are stored as text in this workbook (e.g. "1.00",
# "63.688.01" using '.' as a thousands separator) even though they look
# numeric. A leading apostrophe forces the COM layer to keep them as text
# instead of silently coercing to a float (which would mangle values like
# "84.80" -> 84.8 or "0.0400" -> 0.04, and corrupt the PEPE subscript price
# entirely). The subscript digit itself must be built via [char] and coerced
# to [string] before concatenation, otherwise '+' on a [char] operand does
# numeric addition instead of string concatenation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price($row, $price) {
    $ws.Cells.Item($row, 4).Value = "'" + $price
}

function Set-Volume($row, $volume) {
    $ws.Cells.Item($row, 5).Value = $volume
}

Set-Price  2 "63.754.11"
Set-Volume 2 "  +1.43%  "

Set-Price  3 "3.317.88"
Set-Volume 3 "  +4.95%  "

Set-Volume 4 "  +0.11%  "

Set-Price  5 "604.15"
Set-Volume 5 "  +2.62%  "

Set-Price  6 "142.53"
Set-Volume 6 "  +2.87%  "

Set-Volume 7 "  +0.01%  "

Set-Price  8 "3.315.28"
Set-Volume 8 "  +4.90%  "

Set-Volume 9 "  +0.66%  "

Set-Volume 10 "  +2.73%  "

Set-Price  11 "5.54"
Set-Volume 11 "  +4.60%  "

Set-Price  12 "0.469"
Set-Volume 12 "  +2.30%  "

Set-Volume 13 "  +1.12%  "

Set-Price  14 "34.86"
Set-Volume 14 "  +2.18%  "

Set-Price  15 "3.861.72"
Set-Volume 15 "  +4.94%  "

Set-Volume 16 "  +0.12%  "

Set-Price  17 "3.319.37"
Set-Volume 17 "  +5.08%  "

Set-Price  18 "63.815.20"
Set-Volume 18 "  +1.51%  "

Set-Price  19 "6.88"
Set-Volume 19 "  +3.21%  "

Set-Price  20 "481.12"
Set-Volume 20 "  +1.88%  "

Set-Price  22 "0.737"
Set-Volume 22 "  +4.97%  "

Set-Price  23 "8.09"
Set-Volume 23 "  +4.45%  "

Set-Price  24 "13.72"
Set-Volume 24 "  +5.88%  "

Set-Price  25 "84.80"
Set-Volume 25 "  +1.38%  "

Set-Volume 26 "  -0.02%  "

Set-Price  27 "2.78"
Set-Volume 27 "  +2.60%  "

Set-Volume 28 "  +0.03%  "

Set-Price  29 "7.23"
Set-Volume 29 "  +1.98%  "

Set-Price  30 "8.20"
Set-Volume 30 "  +3.03%  "

Set-Price  31 "2.16"
Set-Volume 31 "  +3.51%  "

Set-Price  32 "28.96"
Set-Volume 32 "  +8.46%  "

Set-Volume 33 "  +0.07%  "

Set-Volume 34 "  +0.37%  "

Set-Volume 35 "  +3.87%  "

Set-Price  36 "6.08"
Set-Volume 36 "  +5.16%  "

Set-Price  37 "52.28"
Set-Volume 37 "  -0.42%  "

$subscript3 = [string][char]0x2083
$pepePrice = "0.0" + $subscript3 + "0741"
Set-Price  38 $pepePrice
Set-Volume 38 "  +5.32%  "

# Rows 39/40 swap: VeChain moves up to row 39 (with refreshed numbers),
# Bittensor moves down to row 40 (with refreshed numbers).
$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-Price  39 "0.0400"
Set-Volume 39 "  +3.27%  "

$ws.Cells.Item(40, 2).Value = "Bittensor"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-Price  40 "436.13"
Set-Volume 40 "  +4.39%  "

Set-Price  41 "3.098.09"
Set-Volume 41 "  +4.72%  "

Set-Price  42 "0.121"
Set-Volume 42 "  +9.31%  "

Set-Price  43 "2.77"
Set-Volume 43 "  +1.08%  "

Set-Volume 44 "  +0.61%  "

Set-Price  45 "0.265"
Set-Volume 45 "  +0.77%  "

Set-Price  46 "2.24"
Set-Volume 46 "  +5.30%  "

Set-Price  47 "37.46"
Set-Volume 47 "  +17.30%  "

Set-Price  48 "26.34"
Set-Volume 48 "  +3.29%  "

Set-Volume 49 "  -0.05%  "

# Rows 50/51 swap: ThetaToken moves up to row 50, Stellar moves down to
# row 51 (prices unchanged, only volumes refreshed).
$ws.Cells.Item(50, 2).Value = "ThetaToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-Price  50 "2.30"
Set-Volume 50 "  +2.46%  "

$ws.Cells.Item(51, 2).Value = "Stellar"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-Price  51 "0.114"
Set-Volume 51 "  +0.51%  "
